# BPW_Lecture_07.pptx — "updated F21 gdrive link"
#
# 1) Refresh the cached "datetimeFigureOut" field text (auto date placeholder)
#    from "1/8/19" to "3/18/21" on the slide master and every slide layout.
# 2) Fix a copy/paste typo on slide 20 ("Rectangle 13"): the 3rd bullet's
#    trailing comment said "forwards then forwards" — it should read
#    "forwards then backwards". The run is split into three runs so only the
#    "then forwards " -> "then backwards " portion actually changes value.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 14) {
            $pf = $sh.PlaceholderFormat
            if ($pf.Type -eq 16) {
                $sh.TextFrame.TextRange.Text = "3/18/21"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# --- slide 20 comment fix -------------------------------------------------

$slide = $p.Slides.Item(20)
$shape = $slide.Shapes.Item("Rectangle 13")
$tr = $shape.TextFrame.TextRange
$para = $tr.Paragraphs(3, 1)
$paraText = $para.Text

$oldComment = "  /* animation plays forwards then forwards */"
$commentIdx0 = $paraText.IndexOf($oldComment)

if ($commentIdx0 -ge 0) {
    $oldPrefix = "  /* animation plays forwards "
    $oldMiddle = "then forwards "
    $oldSuffix = "*/"

    $newPrefix = "  /* animation plays forwards "
    $newMiddle = "then backwards "
    $newSuffix = "*/"

    $commentStart = $commentIdx0 + 1
    $prefixStart = $commentStart
    $middleStart = $prefixStart + $oldPrefix.Length
    $suffixStart = $middleStart + $oldMiddle.Length

    # Edit right-to-left so earlier offsets stay valid as lengths change.
    $suffixRange = $para.Characters($suffixStart, $oldSuffix.Length)
    $suffixRange.Text = $newSuffix

    $middleRange = $para.Characters($middleStart, $oldMiddle.Length)
    $middleRange.Text = $newMiddle

    $prefixRange = $para.Characters($prefixStart, $oldPrefix.Length)
    $prefixRange.Text = $newPrefix
}
